$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add "x" marks in column G for rows 5, 6, 16-21 (matching the style of existing G2:G4 "x" marks)
$ws.Range("G5").Value = "x"
$ws.Range("G6").Value = "x"
$ws.Range("G16").Value = "x"
$ws.Range("G17").Value = "x"
$ws.Range("G18").Value = "x"
$ws.Range("G19").Value = "x"
$ws.Range("G20").Value = "x"
$ws.Range("G21").Value = "x"

# Copy the style (center aligned) from an existing "x" cell (G2) to the new cells
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G5:G6").PasteSpecial(-4122) | Out-Null
$ws.Range("G16:G21").PasteSpecial(-4122) | Out-Null

# Move selection to G22 as the final active cell
$ws.Range("G22").Select() | Out-Null
